$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Q1").EntireColumn.Insert()
$ws.Range("R1").EntireColumn.Insert()
$ws.Range("P1").EntireColumn.ColumnWidth = 15.921875
$ws.Range("Q1").EntireColumn.ColumnWidth = 13.921875
$ws.Range("R1").EntireColumn.ColumnWidth = 15.921875
$ws.Range("S1").EntireColumn.ColumnWidth = 13.921875
